$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell holds plain text (prices/percentages formatted as strings in the
# source data, e.g. "0.470" or "61.475.74"). Assigning a numeric-looking string to a
# "General"-formatted cell would make Excel auto-convert it to a number and lose
# formatting (trailing zeros, thousand-dot grouping, etc.), so we briefly switch the
# cell to Text format, write the literal string, then restore the original format.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.489.47'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.20%  '
$ws.Range("E2").NumberFormat = "General"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.361.89'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.28%  '
$ws.Range("E3").NumberFormat = "General"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("E4").NumberFormat = "General"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.50'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("E5").NumberFormat = "General"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.38'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.60%  '
$ws.Range("E6").NumberFormat = "General"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E7").NumberFormat = "General"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.358.95'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.37%  '
$ws.Range("E8").NumberFormat = "General"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.469'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.56%  '
$ws.Range("E9").NumberFormat = "General"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.47'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.20%  '
$ws.Range("E10").NumberFormat = "General"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.42%  '
$ws.Range("E11").NumberFormat = "General"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.49%  '
$ws.Range("E12").NumberFormat = "General"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.935.40'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.30%  '
$ws.Range("E13").NumberFormat = "General"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.123'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.68%  '
$ws.Range("E14").NumberFormat = "General"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000174'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.76%  '
$ws.Range("E15").NumberFormat = "General"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.86'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.38%  '
$ws.Range("E16").NumberFormat = "General"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.365.12'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.27%  '
$ws.Range("E17").NumberFormat = "General"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.523.91'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("E18").NumberFormat = "General"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.03%  '
$ws.Range("E19").NumberFormat = "General"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.85'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -1.01%  '
$ws.Range("E20").NumberFormat = "General"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.30'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.35%  '
$ws.Range("E21").NumberFormat = "General"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '379.63'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.07%  '
$ws.Range("E22").NumberFormat = "General"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.553'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.77%  '
$ws.Range("E23").NumberFormat = "General"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.499.84'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.26%  '
$ws.Range("E24").NumberFormat = "General"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("E25").NumberFormat = "General"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.23%  '
$ws.Range("E26").NumberFormat = "General"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '71.24'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.29%  '
$ws.Range("E27").NumberFormat = "General"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.79'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +10.70%  '
$ws.Range("E28").NumberFormat = "General"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.50'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.73%  '
$ws.Range("E29").NumberFormat = "General"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.14%  '
$ws.Range("E30").NumberFormat = "General"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.165'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.19%  '
$ws.Range("E31").NumberFormat = "General"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.17'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.98%  '
$ws.Range("E32").NumberFormat = "General"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.67%  '
$ws.Range("E33").NumberFormat = "General"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("E34").NumberFormat = "General"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.56'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.18%  '
$ws.Range("E35").NumberFormat = "General"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.20'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -5.95%  '
$ws.Range("E36").NumberFormat = "General"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.79'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.19%  '
$ws.Range("E37").NumberFormat = "General"
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'Monero'
$ws.Range("B38").NumberFormat = "General"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("C38").NumberFormat = "General"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '165.64'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.31%  '
$ws.Range("E38").NumberFormat = "General"
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'ImmutableX'
$ws.Range("B39").NumberFormat = "General"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("C39").NumberFormat = "General"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.52'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.01%  '
$ws.Range("E39").NumberFormat = "General"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0762'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.03%  '
$ws.Range("E40").NumberFormat = "General"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("E41").NumberFormat = "General"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.12%  '
$ws.Range("E42").NumberFormat = "General"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.67%  '
$ws.Range("E43").NumberFormat = "General"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.88%  '
$ws.Range("E44").NumberFormat = "General"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '41.44'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.90%  '
$ws.Range("E45").NumberFormat = "General"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.38'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.30%  '
$ws.Range("E46").NumberFormat = "General"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '23.91'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("E47").NumberFormat = "General"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.84'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.14%  '
$ws.Range("E48").NumberFormat = "General"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '23.00'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.10%  '
$ws.Range("E49").NumberFormat = "General"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.366.53'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.29%  '
$ws.Range("E50").NumberFormat = "General"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.65%  '
$ws.Range("E51").NumberFormat = "General"
